$wb = $excel.ActiveWorkbook

# --- 1. Refresh the handoff status / timestamps for the remaining
#        (1117757f-...) entry on every sheet, then drop the whole
#        (ad6ac1e5-...) row (row 3) which is no longer part of the report.

# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-48-13 04:48:41"

# Sheet "zh-cn": refresh the shared Status cell + the Latest Handoff Datetime (column E)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-13 04:48:37"

# Sheet "de-de": refresh the shared Status cell + the Latest Handoff Datetime (column E)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-13 04:48:41"

# --- 2. Remove row 3 (the ad6ac1e5-... entry) from all three sheets,
#        including its hyperlinks, since it has been handed off already.
foreach ($ws in @($wsOverview, $wsZhCn, $wsDeDe)) {
    $found = $true
    while ($found) {
        $found = $false
        foreach ($hl in $ws.Hyperlinks) {
            $addr = $hl.Range.Address()
            if ($addr -like "*3") {
                $hl.Delete()
                $found = $true
                break
            }
        }
    }
    $ws.Rows.Item(3).Delete()
}
